# Auto-generated edit script applying the Malboro_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 654.36365
$ws.Cells.Item(12, 9).Value = 683.3333
$ws.Cells.Item(12, 10).Value = 619.6
$ws.Cells.Item(12, 11).Value = 683.3333
$ws.Cells.Item(12, 12).Value = 619.6
$ws.Cells.Item(12, 13).Value = -513.3333
$ws.Cells.Item(12, 14).Value = -959.6
$ws.Cells.Item(40, 8).Value = 892
$ws.Cells.Item(40, 9).Value = 893
$ws.Cells.Item(40, 10).Value = 890
$ws.Cells.Item(40, 11).Value = 893
$ws.Cells.Item(40, 12).Value = 890
$ws.Cells.Item(40, 13).Value = -718
$ws.Cells.Item(40, 14).Value = -1240
$ws.Cells.Item(69, 8).Value = 1000
$ws.Cells.Item(69, 9).Value = 1000
$ws.Cells.Item(69, 11).Value = 3000
$ws.Cells.Item(69, 13).Value = -2126
$ws.Cells.Item(72, 8).Value = 1000
$ws.Cells.Item(72, 9).Value = 1000
$ws.Cells.Item(72, 11).Value = 9000
$ws.Cells.Item(72, 13).Value = -4632
$ws.Cells.Item(112, 8).Value = 7227.2856
$ws.Cells.Item(112, 10).Value = 5769.516
$ws.Cells.Item(112, 12).Value = 17308.548
$ws.Cells.Item(112, 14).Value = -19524.548
$ws.Cells.Item(116, 8).Value = 4275.603
$ws.Cells.Item(116, 9).Value = 3987.3438
$ws.Cells.Item(116, 11).Value = 3987.3438
$ws.Cells.Item(116, 13).Value = -545.3438000000001
$ws.Cells.Item(125, 8).Value = 5347.8823
$ws.Cells.Item(125, 9).Value = 4433.6
$ws.Cells.Item(125, 10).Value = 6654
$ws.Cells.Item(125, 11).Value = 39902.4
$ws.Cells.Item(125, 12).Value = 59886
$ws.Cells.Item(125, 13).Value = -37442.4
$ws.Cells.Item(125, 14).Value = -64806
$ws.Cells.Item(137, 8).Value = 10710.895
$ws.Cells.Item(137, 9).Value = 4556.3335
$ws.Cells.Item(137, 10).Value = 15186.939
$ws.Cells.Item(137, 11).Value = 13669.0005
$ws.Cells.Item(137, 12).Value = 45560.817
$ws.Cells.Item(137, 13).Value = -11119.0005
$ws.Cells.Item(137, 14).Value = -50660.817
$ws.Cells.Item(138, 8).Value = 5612.96
$ws.Cells.Item(138, 9).Value = 4626.8335
$ws.Cells.Item(138, 10).Value = 5924.3687
$ws.Cells.Item(138, 11).Value = 13880.5005
$ws.Cells.Item(138, 12).Value = 17773.1061
$ws.Cells.Item(138, 13).Value = -8740.500499999998
$ws.Cells.Item(138, 14).Value = -28053.1061

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1083567.8
$ws.Cells.Item(32, 9).Value = 1647478.5
$ws.Cells.Item(32, 11).Value = 1647478.5
$ws.Cells.Item(32, 13).Value = -1647191.5
$ws.Cells.Item(74, 8).Value = 10132.418
$ws.Cells.Item(74, 9).Value = 2500.359
$ws.Cells.Item(74, 10).Value = 28735.562
$ws.Cells.Item(74, 11).Value = 2500.359
$ws.Cells.Item(74, 12).Value = 28735.562
$ws.Cells.Item(74, 13).Value = -1626.359
$ws.Cells.Item(74, 14).Value = -30483.562
$ws.Cells.Item(77, 8).Value = 10132.418
$ws.Cells.Item(77, 9).Value = 2500.359
$ws.Cells.Item(77, 10).Value = 28735.562
$ws.Cells.Item(77, 11).Value = 12501.795
$ws.Cells.Item(77, 12).Value = 143677.81
$ws.Cells.Item(77, 13).Value = -8133.795
$ws.Cells.Item(77, 14).Value = -152413.81
$ws.Cells.Item(88, 8).Value = 4003.5
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 4003.5
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 4003.5
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(88, 14).Value = -4815.5
$ws.Cells.Item(91, 8).Value = 4003.5
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 4003.5
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 4003.5
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(91, 14).Value = -6811.5
$ws.Cells.Item(110, 8).Value = 1889.4667
$ws.Cells.Item(110, 9).Value = 1917.2858
$ws.Cells.Item(110, 10).Value = 1500
$ws.Cells.Item(110, 11).Value = 1917.2858
$ws.Cells.Item(110, 12).Value = 1500
$ws.Cells.Item(110, 13).Value = 127.7141999999999
$ws.Cells.Item(110, 14).Value = -5590
$ws.Cells.Item(132, 8).Value = 7427.769
$ws.Cells.Item(132, 9).Value = 2794.4722
$ws.Cells.Item(132, 11).Value = 8383.4166
$ws.Cells.Item(132, 13).Value = -5853.4166

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 9742.737999999999
$ws.Cells.Item(134, 9).Value = 4695.393
$ws.Cells.Item(134, 10).Value = 19837.428
$ws.Cells.Item(134, 11).Value = 14086.179
$ws.Cells.Item(134, 12).Value = 59512.284
$ws.Cells.Item(134, 13).Value = -11551.179
$ws.Cells.Item(134, 14).Value = -64582.284

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 487.02777
$ws.Cells.Item(7, 9).Value = 538.13794
$ws.Cells.Item(7, 11).Value = 538.13794
$ws.Cells.Item(7, 13).Value = -425.13794
$ws.Cells.Item(31, 8).Value = 20898.889
$ws.Cells.Item(31, 9).Value = 6542.72
$ws.Cells.Item(31, 11).Value = 6542.72
$ws.Cells.Item(31, 13).Value = -6247.72
$ws.Cells.Item(34, 8).Value = 20898.889
$ws.Cells.Item(34, 9).Value = 6542.72
$ws.Cells.Item(34, 11).Value = 6542.72
$ws.Cells.Item(34, 13).Value = -6340.72
$ws.Cells.Item(50, 8).Value = 41665
$ws.Cells.Item(50, 9).Value = 39998
$ws.Cells.Item(50, 11).Value = 39998
$ws.Cells.Item(50, 13).Value = -39373
$ws.Cells.Item(58, 8).Value = 11590.387
$ws.Cells.Item(58, 9).Value = 4997.3076
$ws.Cells.Item(58, 11).Value = 4997.3076
$ws.Cells.Item(58, 13).Value = -4794.3076
$ws.Cells.Item(74, 8).Value = 37500
$ws.Cells.Item(74, 10).Value = 45000
$ws.Cells.Item(74, 12).Value = 45000
$ws.Cells.Item(74, 14).Value = -46748
$ws.Cells.Item(77, 8).Value = 37500
$ws.Cells.Item(77, 10).Value = 45000
$ws.Cells.Item(77, 12).Value = 135000
$ws.Cells.Item(77, 14).Value = -143736
$ws.Cells.Item(105, 8).Value = 11057.728
$ws.Cells.Item(105, 9).Value = 16137.143
$ws.Cells.Item(105, 11).Value = 16137.143
$ws.Cells.Item(105, 13).Value = -14390.143
$ws.Cells.Item(136, 8).Value = 11590.387
$ws.Cells.Item(136, 9).Value = 4997.3076
$ws.Cells.Item(136, 11).Value = 14991.9228
$ws.Cells.Item(136, 13).Value = -12441.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1848
$ws.Cells.Item(5, 9).Value = 487.42856
$ws.Cells.Item(5, 10).Value = 4569.143
$ws.Cells.Item(5, 11).Value = 1462.28568
$ws.Cells.Item(5, 12).Value = 13707.429
$ws.Cells.Item(5, 13).Value = -1350.28568
$ws.Cells.Item(5, 14).Value = -13931.429
$ws.Cells.Item(39, 8).Value = 7067.3335
$ws.Cells.Item(39, 10).Value = 9200
$ws.Cells.Item(39, 12).Value = 27600
$ws.Cells.Item(39, 14).Value = -28188
$ws.Cells.Item(55, 8).Value = 2999.75
$ws.Cells.Item(55, 9).Value = 5399.5
$ws.Cells.Item(55, 10).Value = 600
$ws.Cells.Item(55, 11).Value = 16198.5
$ws.Cells.Item(55, 12).Value = 1800
$ws.Cells.Item(55, 13).Value = -16021.5
$ws.Cells.Item(55, 14).Value = -2154
$ws.Cells.Item(101, 8).Value = 12214.5
$ws.Cells.Item(101, 10).Value = 12214.5
$ws.Cells.Item(101, 12).Value = 36643.5
$ws.Cells.Item(101, 14).Value = -41511.5
$ws.Cells.Item(109, 8).Value = 3335623.8
$ws.Cells.Item(109, 9).Value = 2141
$ws.Cells.Item(109, 10).Value = 6669106.5
$ws.Cells.Item(109, 11).Value = 6423
$ws.Cells.Item(109, 12).Value = 20007319.5
$ws.Cells.Item(109, 13).Value = -5383
$ws.Cells.Item(109, 14).Value = -20009399.5
$ws.Cells.Item(131, 8).Value = 1482.78
$ws.Cells.Item(131, 10).Value = 1482.78
$ws.Cells.Item(131, 12).Value = 4448.34
$ws.Cells.Item(131, 14).Value = -14528.34
$ws.Cells.Item(135, 8).Value = 1848
$ws.Cells.Item(135, 9).Value = 487.42856
$ws.Cells.Item(135, 10).Value = 4569.143
$ws.Cells.Item(135, 11).Value = 4386.85704
$ws.Cells.Item(135, 12).Value = 41122.287
$ws.Cells.Item(135, 13).Value = -1851.85704
$ws.Cells.Item(135, 14).Value = -46192.287
$ws.Cells.Item(139, 8).Value = 6281.55
$ws.Cells.Item(139, 9).Value = 5588.4546
$ws.Cells.Item(139, 11).Value = 16765.3638
$ws.Cells.Item(139, 13).Value = -11625.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 16128.883
$ws.Cells.Item(39, 10).Value = 16128.883
$ws.Cells.Item(39, 12).Value = 16128.883
$ws.Cells.Item(39, 14).Value = -17192.883
$ws.Cells.Item(52, 8).Value = 47596
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 47596
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 47596
$ws.Cells.Item(52, 13).ClearContents()
$ws.Cells.Item(52, 14).Value = -48114
$ws.Cells.Item(132, 8).Value = 7613.5
$ws.Cells.Item(132, 9).Value = 4603.3706
$ws.Cells.Item(132, 10).Value = 34704.668
$ws.Cells.Item(132, 11).Value = 13810.1118
$ws.Cells.Item(132, 12).Value = 104114.004
$ws.Cells.Item(132, 13).Value = -11280.1118
$ws.Cells.Item(132, 14).Value = -109174.004
$ws.Cells.Item(134, 8).Value = 62142.43
$ws.Cells.Item(134, 10).Value = 62142.43
$ws.Cells.Item(134, 12).Value = 186427.29
$ws.Cells.Item(134, 14).Value = -191497.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3638.7856
$ws.Cells.Item(46, 10).Value = 3272.0908
$ws.Cells.Item(46, 12).Value = 3272.0908
$ws.Cells.Item(46, 14).Value = -3648.0908
$ws.Cells.Item(132, 8).Value = 5203.3706
$ws.Cells.Item(132, 9).Value = 2749.4092
$ws.Cells.Item(132, 10).Value = 16000.8
$ws.Cells.Item(132, 11).Value = 8248.2276
$ws.Cells.Item(132, 12).Value = 48002.39999999999
$ws.Cells.Item(132, 13).Value = -5718.2276
$ws.Cells.Item(132, 14).Value = -53062.39999999999
$ws.Cells.Item(136, 8).Value = 14898.082
$ws.Cells.Item(136, 9).Value = 12757.849
$ws.Cells.Item(136, 10).Value = 19312.312
$ws.Cells.Item(136, 11).Value = 38273.547
$ws.Cells.Item(136, 12).Value = 57936.936
$ws.Cells.Item(136, 13).Value = -35723.547
$ws.Cells.Item(136, 14).Value = -63036.936

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 9287.793
$ws.Cells.Item(136, 9).Value = 1875.7391
$ws.Cells.Item(136, 10).Value = 37700.668
$ws.Cells.Item(136, 11).Value = 5627.2173
$ws.Cells.Item(136, 12).Value = 113102.004
$ws.Cells.Item(136, 13).Value = -3077.2173
$ws.Cells.Item(136, 14).Value = -118202.004
